$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.464.43'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '1.858.75'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.95'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4772'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3792'
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07308'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9289'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.71'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07785'
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').Value = '1.867.23'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.451'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.22'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.012'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008820'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '27.486.49'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.61'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.095'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.939'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.79'
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.44'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.001'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.21'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.934'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08888'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.332'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.203'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7512'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.577'
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.715'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02043'
$ws.Range('E36').Value = '  +4.29%  '
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5563'
$ws.Range('E38').Value = '  +5.84%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.988'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.020'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.569'
$ws.Range('E42').Value = '  +3.47%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4866'
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('E47').Value = '  +3.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.12'
$ws.Range('E48').Value = '  +1.28%  '
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9135'
$ws.Range('E51').Value = '  +2.78%  '
